$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.11327129920958612,
    -0.0059999999397959414,
    0.028223395153343489,
    -0.0079999999079340967,
    -0.0029999999527037247,
    -0.0019999999513693467,
    -0.0099999998704403126,
    -0.0099999998650446287,
    -0.0019999999371358435,
    -0.0019999999304616267,
    -0.002999999919755858,
    -0.0034999999133540349,
    -0.0034999999068610066,
    -0.0079999998607727107,
    -0.00099999992842114693,
    -0.0019999999174102889,
    -0.001999999915717865,
    -0.0039999998955986271,
    -0.0039999999592792435,
    -0.0039999999561164401,
    -0.003999999955652811,
    -0.0039999999552522425,
    -0.0049999999382794869,
    -0.01999999978470246,
    -0.019999999781670219,
    -0.0024999999462398392,
    0.018884908676248369,
    -0.0019999999461557394,
    -0.0069999998939733743,
    -0.059999999372432278,
    0.062499643016703743,
    -0.0099999998722370975,
    -0.0039999999310307288
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

